$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -830
$ws.Range("H17").Value = 874.6957
$ws.Range("J17").Value = 807.13635
$ws.Range("L17").Value = 2421.40905
$ws.Range("N17").Value = -2757.40905
$ws.Range("H32").Value = 20004900
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 20004900
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 20004900
$ws.Range("M32").Value = $null
$ws.Range("N32").Value = -20005552
$ws.Range("H33").Value = 525.13336
$ws.Range("I33").Value = 334.36365
$ws.Range("J33").Value = 1049.75
$ws.Range("K33").Value = 334.36365
$ws.Range("L33").Value = 1049.75
$ws.Range("M33").Value = -105.36365
$ws.Range("N33").Value = -1507.75
$ws.Range("H40").Value = 4126.4243
$ws.Range("I40").Value = 3432.5
$ws.Range("J40").Value = 4637.737
$ws.Range("K40").Value = 3432.5
$ws.Range("L40").Value = 4637.737
$ws.Range("M40").Value = -3257.5
$ws.Range("N40").Value = -4987.737
$ws.Range("H41").Value = 6859.8125
$ws.Range("I41").Value = 996.3333
$ws.Range("J41").Value = 10377.9
$ws.Range("K41").Value = 996.3333
$ws.Range("L41").Value = 10377.9
$ws.Range("M41").Value = -556.3333
$ws.Range("N41").Value = -11257.9
$ws.Range("H99").Value = 7661.75
$ws.Range("I99").Value = 184.85715
$ws.Range("K99").Value = 554.5714499999999
$ws.Range("M99").Value = 943.4285500000001
$ws.Range("H118").Value = 1446.6666
$ws.Range("I118").Value = 1431.4286
$ws.Range("K118").Value = 4294.2858
$ws.Range("M118").Value = -2637.2858
$ws.Range("H137").Value = 2508.4565
$ws.Range("I137").Value = 2110.0732
$ws.Range("J137").Value = 5775.2
$ws.Range("K137").Value = 6330.219599999999
$ws.Range("L137").Value = 17325.6
$ws.Range("M137").Value = -3780.219599999999
$ws.Range("N137").Value = -22425.6
$ws.Range("H138").Value = 2817.3284
$ws.Range("J138").Value = 3629.6191
$ws.Range("L138").Value = 10888.8573
$ws.Range("N138").Value = -21168.8573

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 3337.6667
$ws.Range("I19").Value = 1989
$ws.Range("J19").Value = 4012
$ws.Range("K19").Value = 1989
$ws.Range("L19").Value = 4012
$ws.Range("M19").Value = -1760
$ws.Range("N19").Value = -4470
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = $null
$ws.Range("H32").Value = 22731380
$ws.Range("I32").Value = 24394396
$ws.Range("K32").Value = 24394396
$ws.Range("M32").Value = -24394109
$ws.Range("H132").Value = 1815.5758
$ws.Range("I132").Value = 1815.5758
$ws.Range("K132").Value = 5446.7274
$ws.Range("M132").Value = -2916.7274

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2003.3914
$ws.Range("I20").Value = 1321.2667
$ws.Range("K20").Value = 1321.2667
$ws.Range("M20").Value = -1074.2667
$ws.Range("H30").Value = 1500
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").Value = $null
$ws.Range("H86").Value = 4435
$ws.Range("I86").Value = 2879.8
$ws.Range("J86").Value = 5545.857
$ws.Range("K86").Value = 2879.8
$ws.Range("L86").Value = 5545.857
$ws.Range("M86").Value = -1756.8
$ws.Range("N86").Value = -7791.857
$ws.Range("H89").Value = 4435
$ws.Range("I89").Value = 2879.8
$ws.Range("J89").Value = 5545.857
$ws.Range("K89").Value = 14399
$ws.Range("L89").Value = 27729.285
$ws.Range("M89").Value = -8783
$ws.Range("N89").Value = -38961.285
$ws.Range("H134").Value = 2448.3845
$ws.Range("I134").Value = 1859.0435
$ws.Range("K134").Value = 5577.1305
$ws.Range("M134").Value = -3042.1305

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6490.737
$ws.Range("I134").Value = 6572.778
$ws.Range("K134").Value = 19718.334
$ws.Range("M134").Value = -17183.334

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1585.5294
$ws.Range("I5").Value = 1561.3572
$ws.Range("J5").Value = 1698.3334
$ws.Range("K5").Value = 4684.071599999999
$ws.Range("L5").Value = 5095.0002
$ws.Range("M5").Value = -4572.071599999999
$ws.Range("N5").Value = -5319.0002
$ws.Range("H8").Value = 1668.4445
$ws.Range("I8").Value = 1668.4445
$ws.Range("K8").Value = 5005.333500000001
$ws.Range("M8").Value = -4866.333500000001
$ws.Range("H10").Value = 178.66667
$ws.Range("I10").Value = 63.5
$ws.Range("K10").Value = 190.5
$ws.Range("M10").Value = -51.5
$ws.Range("H12").Value = 2.7692308
$ws.Range("J12").Value = 2.625
$ws.Range("L12").Value = 7.875
$ws.Range("N12").Value = -353.875
$ws.Range("H75").Value = 12650.2
$ws.Range("J75").Value = 15212.5
$ws.Range("L75").Value = 45637.5
$ws.Range("N75").Value = -47633.5
$ws.Range("H78").Value = 12650.2
$ws.Range("J78").Value = 15212.5
$ws.Range("L78").Value = 136912.5
$ws.Range("N78").Value = -146896.5
$ws.Range("H107").Value = 2273.625
$ws.Range("I107").Value = 1842.2222
$ws.Range("J107").Value = 2828.2856
$ws.Range("K107").Value = 5526.6666
$ws.Range("L107").Value = 8484.856800000001
$ws.Range("M107").Value = -3606.6666
$ws.Range("N107").Value = -12324.8568
$ws.Range("H113").Value = 1566.2727
$ws.Range("J113").Value = 1797.4286
$ws.Range("L113").Value = 5392.2858
$ws.Range("N113").Value = -9732.2858
$ws.Range("H121").Value = 2053.6365
$ws.Range("I121").Value = 850
$ws.Range("K121").Value = 2550
$ws.Range("M121").Value = -1240
$ws.Range("H133").Value = 2686.3333
$ws.Range("I133").Value = 2686.3333
$ws.Range("K133").Value = 8058.999899999999
$ws.Range("M133").Value = -2998.999899999999
$ws.Range("H135").Value = 1585.5294
$ws.Range("I135").Value = 1561.3572
$ws.Range("J135").Value = 1698.3334
$ws.Range("K135").Value = 14052.2148
$ws.Range("L135").Value = 15285.0006
$ws.Range("M135").Value = -11517.2148
$ws.Range("N135").Value = -20355.0006
$ws.Range("H138").Value = 4965.5
$ws.Range("I138").Value = 4418.933
$ws.Range("K138").Value = 13256.799
$ws.Range("M138").Value = -8116.798999999999

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 6273.4
$ws.Range("I36").Value = 955.6667
$ws.Range("K36").Value = 955.6667
$ws.Range("M36").Value = -470.6667
$ws.Range("H122").Value = 4696.476
$ws.Range("I122").Value = 3821.5715
$ws.Range("J122").Value = 6446.2856
$ws.Range("K122").Value = 11464.7145
$ws.Range("L122").Value = 19338.8568
$ws.Range("M122").Value = -9014.7145
$ws.Range("N122").Value = -24238.8568
$ws.Range("H132").Value = 4268.7393
$ws.Range("I132").Value = 4509.1
$ws.Range("J132").Value = 2666.3333
$ws.Range("K132").Value = 13527.3
$ws.Range("L132").Value = 7998.999899999999
$ws.Range("M132").Value = -10997.3
$ws.Range("N132").Value = -13058.9999

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5348.643
$ws.Range("I7").Value = 4984.3335
$ws.Range("K7").Value = 4984.3335
$ws.Range("M7").Value = -4872.3335
$ws.Range("H126").Value = 5348.643
$ws.Range("I126").Value = 4984.3335
$ws.Range("K126").Value = 14953.0005
$ws.Range("M126").Value = -12483.0005
$ws.Range("H136").Value = 2689.3635
$ws.Range("I136").Value = 2448.3928
$ws.Range("K136").Value = 7345.178400000001
$ws.Range("M136").Value = -4795.178400000001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13775.1
$ws.Range("J45").Value = 15226.875
$ws.Range("L45").Value = 15226.875
$ws.Range("N45").Value = -16208.875
